# Mise à jour de l'application
# Appends the 2026-01-28 (serial 46050) wellness entries for every player
# after the last existing row (773), extending the data table through
# row 787, and updates the sheet selection to reflect where entry stopped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New rows to append (row, player name, volume, intensite, fatigue,
# douleur, localisation douleur [or $null when empty], plaisir).
# Charge (col I) = Volume * Intensite, same formula as the rest of the
# column.
# ---------------------------------------------------------------------
# Non-breaking space: a couple of the existing "Localisation douleur"
# shared strings end with U+00A0 rather than a plain space, so the new
# rows that reuse those exact labels need the same character to land on
# the pre-existing shared-string entry instead of minting a duplicate.
$nbsp = [char]0x00A0

$newRows = @(
    @{ Row=774; Nom="Amir Etien";        Volume=70; Intensite=5; Fatigue=7; Douleur=5; Lieu="Ischio";                 Plaisir=2 },
    @{ Row=775; Nom="Yoann Martelat";    Volume=70; Intensite=5; Fatigue=6; Douleur=4; Lieu="Genou";                  Plaisir=5 },
    @{ Row=776; Nom="Romain Thunet";     Volume=70; Intensite=5; Fatigue=5; Douleur=3; Lieu=("Synthétique" + $nbsp);  Plaisir=5 },
    @{ Row=777; Nom="Naim Ighbane";      Volume=70; Intensite=6; Fatigue=6; Douleur=2; Lieu="Coup tibia";             Plaisir=7 },
    @{ Row=778; Nom="Malik Boussaid";    Volume=70; Intensite=7; Fatigue=7; Douleur=0; Lieu=$null;                    Plaisir=8 },
    @{ Row=779; Nom="Yoan Zouma";        Volume=70; Intensite=2; Fatigue=6; Douleur=0; Lieu=$null;                    Plaisir=5 },
    @{ Row=780; Nom="Rayane Chayebi";    Volume=70; Intensite=4; Fatigue=7; Douleur=2; Lieu="Courbature";             Plaisir=6 },
    @{ Row=781; Nom="Mattheo Haon";      Volume=70; Intensite=4; Fatigue=0; Douleur=0; Lieu=$null;                    Plaisir=3 },
    @{ Row=782; Nom="Ilan Ihaddadene";   Volume=70; Intensite=4; Fatigue=8; Douleur=0; Lieu=$null;                    Plaisir=3 },
    @{ Row=783; Nom="Naim Dhib";         Volume=70; Intensite=6; Fatigue=3; Douleur=3; Lieu="Psoas";                  Plaisir=6 },
    @{ Row=784; Nom="Sofiane Belle";     Volume=70; Intensite=4; Fatigue=4; Douleur=3; Lieu="Talon";                  Plaisir=1 },
    @{ Row=785; Nom="Theo Owono";        Volume=70; Intensite=3; Fatigue=3; Douleur=0; Lieu=$null;                    Plaisir=8 },
    @{ Row=786; Nom="Levy Ndoutoume";    Volume=70; Intensite=5; Fatigue=7; Douleur=9; Lieu=("Adducteur" + $nbsp);    Plaisir=2 },
    @{ Row=787; Nom="Karahali Souaré";   Volume=70; Intensite=5; Fatigue=6; Douleur=6; Lieu="Cheville";               Plaisir=8 }
)

$dateSerial = 46050   # 2026-01-28

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Clone the formatting of the preceding row that has the same shape
    # (row 772 has a filled "Localisation douleur" cell, row 773 has an
    # empty one) so the new row's style indices match exactly.
    if ($null -eq $entry.Lieu) {
        $ws.Range("A773:I773").Copy()
    } else {
        $ws.Range("A772:I772").Copy()
    }
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)

    $ws.Range("A$r").Value = $dateSerial
    $ws.Range("B$r").Value = $entry.Nom
    $ws.Range("C$r").Value = $entry.Volume
    $ws.Range("D$r").Value = $entry.Intensite
    $ws.Range("E$r").Value = $entry.Fatigue
    $ws.Range("F$r").Value = $entry.Douleur
    if ($null -ne $entry.Lieu) {
        $ws.Range("G$r").Value = $entry.Lieu
    }
    $ws.Range("H$r").Value = $entry.Plaisir
    $ws.Range("I$r").Formula = "=C$r*D$r"
}

$excel.CutCopyMode = 0

# Recalculate so every Charge cell carries a fresh cached value.
$excel.Calculate()

# Reflect where data entry left off.
$ws.Range("J772").Select()
